$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New feedback rows to append (all values stored as text, matching the
# existing rows in the sheet, which are all typed as strings).
$newRows = @(
    @("1770957814785", "Alfin Sen Varghese", "alfinsen@gmail.com", "5", "hello..", "2/13/2026, 10:13:34 AM"),
    @("1771848606334", "Alfin", "alfin@123", "5", "good!", "2/23/2026, 5:40:06 PM")
)

$startRow = 6
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        # Force text storage (numeric-looking values like IDs/ratings must
        # not be reinterpreted as numbers), then restore the default
        # "Normal" style so the new cells don't pick up a stray custom
        # number-format style like the rest of the sheet.
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$c - 1]
        $cell.Style = "Normal"
    }
}
